# Updating for version 20.4.1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (pushes protocol/port/os from E,F,G -> F,G,H)
# and shifts in a new "enable_password" column holding a copy of the
# password value for each host.
$ws.Range("E1").EntireColumn.Insert()

# Set the width for the newly inserted column E (target stored width ~17.6640625)
$ws.Range("E1").ColumnWidth = 16.8

# Header for new column
$ws.Range("E1").Value = "enable_password"

# Populate enable_password values = same as password column (D) for each row
$ws.Range("E2").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("D3").Value2

# Username for csr1000v-1 changes from "root" to "developer"
$ws.Range("C3").Value = "developer"

# Update the selected cell to match the saved state
$ws.Range("E7").Select()
